# Add team Wins/Losses/Ties columns (AD/AE/AF) to the MIN_2010 roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new labels, formatted like the existing headers ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the existing header formatting (bold, bordered, centered) from AC1
# onto the new header cells so the new columns visually match row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows (2-44): same W/L/T record repeated for every player row ---
$lastRow = 44
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("AD$r").Value = 94
    $ws.Range("AE$r").Value = 68
    $ws.Range("AF$r").Value = 0
}
